$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poland IV Liga")

# Apply the league-database refresh: row-data corrections (home/away swaps,
# corrected scores/odds/ids) plus the two shared-team-name relabelings
# (Tarnovia Tarnowo Podgorne <-> MGKS Moto Jelcz Olawa, Ursus Warsawa <-> Spartakus Daleszyce).

# Row 17
$ws.Cells.Item(17, 6).Value = 'MGKS Moto Jelcz Olawa'
# Row 45
$ws.Cells.Item(45, 2).Value = 7068599
$ws.Cells.Item(45, 5).Value = 'LKS Jawiszowice'
$ws.Cells.Item(45, 6).Value = 'MKS TrzebiniaSiersza'
$ws.Cells.Item(45, 7).Value = 1
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 11).Value = 'H'
$ws.Cells.Item(45, 12).Value = 1.5
$ws.Cells.Item(45, 13).Value = 4.5
$ws.Cells.Item(45, 14).Value = 4.333
$ws.Cells.Item(45, 15).Value = 1.5
$ws.Cells.Item(45, 16).Value = 4.5
$ws.Cells.Item(45, 17).Value = 4.333
$ws.Cells.Item(45, 18).Value = -1
$ws.Cells.Item(45, 19).Value = 1.8
$ws.Cells.Item(45, 20).Value = 2
$ws.Cells.Item(45, 21).Value = 3.25
$ws.Cells.Item(45, 22).Value = 1.85
$ws.Cells.Item(45, 23).Value = 1.95
$ws.Cells.Item(45, 24).Value = 0.5
$ws.Cells.Item(45, 25).Value = -1
$ws.Cells.Item(45, 27).Value = 0
$ws.Cells.Item(45, 28).Value = 0
$ws.Cells.Item(45, 29).Value = -1
$ws.Cells.Item(45, 30).Value = 0.95
# Row 46
$ws.Cells.Item(46, 5).Value = 'Tarnovia Tarnowo Podgorne'
# Row 47
$ws.Cells.Item(47, 2).Value = 7068602
$ws.Cells.Item(47, 5).Value = 'MGKS Moto Jelcz Olawa'
$ws.Cells.Item(47, 6).Value = 'Victoria Wrzesnia'
$ws.Cells.Item(47, 7).Value = 2
$ws.Cells.Item(47, 8).Value = 2
$ws.Cells.Item(47, 11).Value = 'D'
$ws.Cells.Item(47, 12).Value = 2.1
$ws.Cells.Item(47, 13).Value = 4.2
$ws.Cells.Item(47, 14).Value = 2.5
$ws.Cells.Item(47, 15).Value = 2.1
$ws.Cells.Item(47, 16).Value = 4.2
$ws.Cells.Item(47, 17).Value = 2.55
$ws.Cells.Item(47, 18).Value = -0.25
$ws.Cells.Item(47, 19).Value = 1.95
$ws.Cells.Item(47, 20).Value = 1.85
$ws.Cells.Item(47, 21).Value = 3
$ws.Cells.Item(47, 22).Value = 1.775
$ws.Cells.Item(47, 23).Value = 2.025
$ws.Cells.Item(47, 24).Value = -1
$ws.Cells.Item(47, 25).Value = 3.2
$ws.Cells.Item(47, 27).Value = -0.5
$ws.Cells.Item(47, 28).Value = 0.425
$ws.Cells.Item(47, 29).Value = 0.7749999999999999
$ws.Cells.Item(47, 30).Value = -1
# Row 65
$ws.Cells.Item(65, 2).Value = 7140479
$ws.Cells.Item(65, 6).Value = 'Korona Kielce II'
$ws.Cells.Item(65, 12).Value = 3.1
$ws.Cells.Item(65, 14).Value = 1.85
$ws.Cells.Item(65, 15).Value = 4
$ws.Cells.Item(65, 17).Value = 1.666
$ws.Cells.Item(65, 18).Value = 0.75
$ws.Cells.Item(65, 19).Value = 1.95
$ws.Cells.Item(65, 20).Value = 1.85
$ws.Cells.Item(65, 22).Value = 1.85
$ws.Cells.Item(65, 23).Value = 1.95
$ws.Cells.Item(65, 26).Value = 0.6659999999999999
$ws.Cells.Item(65, 28).Value = 0.8500000000000001
$ws.Cells.Item(65, 29).Value = 0.8500000000000001
# Row 66
$ws.Cells.Item(66, 2).Value = 7140477
$ws.Cells.Item(66, 6).Value = 'Wisla Plock II'
$ws.Cells.Item(66, 12).Value = 2.4
$ws.Cells.Item(66, 14).Value = 2.25
$ws.Cells.Item(66, 15).Value = 2.75
$ws.Cells.Item(66, 17).Value = 2
$ws.Cells.Item(66, 18).Value = 0.25
$ws.Cells.Item(66, 19).Value = 1.975
$ws.Cells.Item(66, 20).Value = 1.825
$ws.Cells.Item(66, 22).Value = 2.025
$ws.Cells.Item(66, 23).Value = 1.775
$ws.Cells.Item(66, 26).Value = 1
$ws.Cells.Item(66, 28).Value = 0.825
$ws.Cells.Item(66, 29).Value = 1.025
# Row 72
$ws.Cells.Item(72, 2).Value = 7183407
$ws.Cells.Item(72, 5).Value = 'Prochowiczanka Prochowice'
$ws.Cells.Item(72, 6).Value = 'Miedz Legnica II'
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 12).Value = 2.25
$ws.Cells.Item(72, 13).Value = 3.75
$ws.Cells.Item(72, 14).Value = 2.5
$ws.Cells.Item(72, 15).Value = 4.5
$ws.Cells.Item(72, 16).Value = 4.2
$ws.Cells.Item(72, 17).Value = 1.533
$ws.Cells.Item(72, 18).Value = 1
$ws.Cells.Item(72, 19).Value = 1.975
$ws.Cells.Item(72, 20).Value = 1.825
$ws.Cells.Item(72, 26).Value = 0.5329999999999999
$ws.Cells.Item(72, 28).Value = 0.825
# Row 73
$ws.Cells.Item(73, 2).Value = 7183410
$ws.Cells.Item(73, 5).Value = 'Arka Gdynia II'
$ws.Cells.Item(73, 6).Value = 'Grom Nowy Staw'
$ws.Cells.Item(73, 7).Value = 1
$ws.Cells.Item(73, 9).Value = 1
$ws.Cells.Item(73, 12).Value = 2.6
$ws.Cells.Item(73, 13).Value = 3.6
$ws.Cells.Item(73, 14).Value = 2.2
$ws.Cells.Item(73, 15).Value = 2.6
$ws.Cells.Item(73, 16).Value = 3.6
$ws.Cells.Item(73, 17).Value = 2.2
$ws.Cells.Item(73, 18).Value = 0
$ws.Cells.Item(73, 19).Value = 2.05
$ws.Cells.Item(73, 20).Value = 1.75
$ws.Cells.Item(73, 26).Value = 1.2
$ws.Cells.Item(73, 28).Value = 0.75
# Row 95
$ws.Cells.Item(95, 6).Value = 'Ursus Warsawa'
# Row 102
$ws.Cells.Item(102, 6).Value = 'MGKS Moto Jelcz Olawa'
# Row 110
$ws.Cells.Item(110, 5).Value = 'Ursus Warsawa'
# Row 123
$ws.Cells.Item(123, 2).Value = 7995829
$ws.Cells.Item(123, 5).Value = 'Korona Kielce II'
$ws.Cells.Item(123, 6).Value = 'Orleta Kielce'
$ws.Cells.Item(123, 7).Value = 2
$ws.Cells.Item(123, 8).Value = 1
$ws.Cells.Item(123, 9).Value = 1
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 11).Value = 'H'
$ws.Cells.Item(123, 12).Value = 1.363
$ws.Cells.Item(123, 13).Value = 5
$ws.Cells.Item(123, 14).Value = 5.5
$ws.Cells.Item(123, 15).Value = 1.222
$ws.Cells.Item(123, 16).Value = 6.5
$ws.Cells.Item(123, 17).Value = 7.5
$ws.Cells.Item(123, 18).Value = -1.75
$ws.Cells.Item(123, 19).Value = 1.8
$ws.Cells.Item(123, 20).Value = 2
$ws.Cells.Item(123, 21).Value = 3.5
$ws.Cells.Item(123, 22).Value = 1.85
$ws.Cells.Item(123, 23).Value = 1.95
$ws.Cells.Item(123, 24).Value = 0.222
$ws.Cells.Item(123, 26).Value = -1
$ws.Cells.Item(123, 28).Value = 1
$ws.Cells.Item(123, 30).Value = 0.95
# Row 124
$ws.Cells.Item(124, 2).Value = 7995827
$ws.Cells.Item(124, 5).Value = 'Alit Ozarow'
$ws.Cells.Item(124, 6).Value = 'GKS Rudki'
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 2
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 1
$ws.Cells.Item(124, 11).Value = 'A'
$ws.Cells.Item(124, 12).Value = 1.909
$ws.Cells.Item(124, 13).Value = 4
$ws.Cells.Item(124, 14).Value = 2.9
$ws.Cells.Item(124, 15).Value = 1.666
$ws.Cells.Item(124, 16).Value = 4.5
$ws.Cells.Item(124, 17).Value = 3.6
$ws.Cells.Item(124, 18).Value = -0.75
$ws.Cells.Item(124, 19).Value = 1.875
$ws.Cells.Item(124, 20).Value = 1.925
$ws.Cells.Item(124, 21).Value = 2.75
$ws.Cells.Item(124, 22).Value = 1.9
$ws.Cells.Item(124, 23).Value = 1.9
$ws.Cells.Item(124, 24).Value = -1
$ws.Cells.Item(124, 26).Value = 2.6
$ws.Cells.Item(124, 28).Value = 0.925
$ws.Cells.Item(124, 30).Value = 0.8999999999999999
# Row 129
$ws.Cells.Item(129, 5).Value = 'Spartakus Daleszyce'
# Row 139
$ws.Cells.Item(139, 6).Value = 'Ursus Warsawa'
# Row 167
$ws.Cells.Item(167, 5).Value = 'Spartakus Daleszyce'
# Row 180
$ws.Cells.Item(180, 6).Value = 'MGKS Moto Jelcz Olawa'
# Row 183
$ws.Cells.Item(183, 2).Value = 8229175
$ws.Cells.Item(183, 5).Value = 'Wierna Malogoszcz'
$ws.Cells.Item(183, 6).Value = 'Klimontowianka Klimontow'
$ws.Cells.Item(183, 8).Value = 3
$ws.Cells.Item(183, 9).Value = 1
$ws.Cells.Item(183, 11).Value = 'A'
$ws.Cells.Item(183, 12).Value = 2.25
$ws.Cells.Item(183, 13).Value = 3.75
$ws.Cells.Item(183, 14).Value = 2.5
$ws.Cells.Item(183, 15).Value = 2.25
$ws.Cells.Item(183, 16).Value = 3.75
$ws.Cells.Item(183, 17).Value = 2.5
$ws.Cells.Item(183, 18).Value = 0
$ws.Cells.Item(183, 19).Value = 1.8
$ws.Cells.Item(183, 20).Value = 2
$ws.Cells.Item(183, 22).Value = 1.8
$ws.Cells.Item(183, 23).Value = 2
$ws.Cells.Item(183, 24).Value = -1
$ws.Cells.Item(183, 26).Value = 1.5
$ws.Cells.Item(183, 27).Value = -1
$ws.Cells.Item(183, 28).Value = 1
$ws.Cells.Item(183, 29).Value = 0.8
$ws.Cells.Item(183, 30).Value = -1
# Row 184
$ws.Cells.Item(184, 2).Value = 8229582
$ws.Cells.Item(184, 5).Value = 'Wisla Krakow II'
$ws.Cells.Item(184, 6).Value = 'Dalin Myslenice'
$ws.Cells.Item(184, 8).Value = 0
$ws.Cells.Item(184, 9).Value = 0
$ws.Cells.Item(184, 11).Value = 'H'
$ws.Cells.Item(184, 12).Value = 1.2
$ws.Cells.Item(184, 13).Value = 6.5
$ws.Cells.Item(184, 14).Value = 8
$ws.Cells.Item(184, 15).Value = 1.285
$ws.Cells.Item(184, 16).Value = 5.5
$ws.Cells.Item(184, 17).Value = 6.25
$ws.Cells.Item(184, 18).Value = -1.75
$ws.Cells.Item(184, 19).Value = 2
$ws.Cells.Item(184, 20).Value = 1.8
$ws.Cells.Item(184, 22).Value = 1.775
$ws.Cells.Item(184, 23).Value = 2.025
$ws.Cells.Item(184, 24).Value = 0.2849999999999999
$ws.Cells.Item(184, 26).Value = -1
$ws.Cells.Item(184, 27).Value = 0.5
$ws.Cells.Item(184, 28).Value = -0.5
$ws.Cells.Item(184, 29).Value = -1
$ws.Cells.Item(184, 30).Value = 1.025
# Row 198
$ws.Cells.Item(198, 6).Value = 'Ursus Warsawa'
# Row 201
$ws.Cells.Item(201, 2).Value = 8257020
$ws.Cells.Item(201, 5).Value = 'GKS Nowiny'
$ws.Cells.Item(201, 6).Value = 'Moravia Morawica'
$ws.Cells.Item(201, 7).Value = 2
$ws.Cells.Item(201, 8).Value = 2
$ws.Cells.Item(201, 10).Value = 2
$ws.Cells.Item(201, 12).Value = 2.2
$ws.Cells.Item(201, 13).Value = 3.75
$ws.Cells.Item(201, 14).Value = 2.55
$ws.Cells.Item(201, 15).Value = 1.9
$ws.Cells.Item(201, 16).Value = 4.1
$ws.Cells.Item(201, 17).Value = 2.9
$ws.Cells.Item(201, 18).Value = -0.5
$ws.Cells.Item(201, 19).Value = 1.975
$ws.Cells.Item(201, 20).Value = 1.825
$ws.Cells.Item(201, 21).Value = 3.75
$ws.Cells.Item(201, 25).Value = 3.1
$ws.Cells.Item(201, 27).Value = -1
$ws.Cells.Item(201, 28).Value = 0.825
$ws.Cells.Item(201, 29).Value = 0.4
$ws.Cells.Item(201, 30).Value = -0.5
# Row 202
$ws.Cells.Item(202, 2).Value = 8256956
$ws.Cells.Item(202, 5).Value = 'Piast Gliwice II'
$ws.Cells.Item(202, 6).Value = 'Rozwoj Katowice'
$ws.Cells.Item(202, 8).Value = 1
$ws.Cells.Item(202, 10).Value = 0
$ws.Cells.Item(202, 11).Value = 'D'
$ws.Cells.Item(202, 12).Value = 2
$ws.Cells.Item(202, 13).Value = 3.5
$ws.Cells.Item(202, 14).Value = 3
$ws.Cells.Item(202, 15).Value = 2
$ws.Cells.Item(202, 16).Value = 3.5
$ws.Cells.Item(202, 17).Value = 3
$ws.Cells.Item(202, 18).Value = -0.25
$ws.Cells.Item(202, 19).Value = 1.775
$ws.Cells.Item(202, 20).Value = 2.025
$ws.Cells.Item(202, 21).Value = 3
$ws.Cells.Item(202, 22).Value = 1.8
$ws.Cells.Item(202, 23).Value = 2
$ws.Cells.Item(202, 25).Value = 2.5
$ws.Cells.Item(202, 26).Value = -1
$ws.Cells.Item(202, 27).Value = -0.5
$ws.Cells.Item(202, 28).Value = 0.5125
$ws.Cells.Item(202, 30).Value = 1
# Row 204
$ws.Cells.Item(204, 2).Value = 8256855
$ws.Cells.Item(204, 5).Value = 'Termalica BB Nieciecza II'
$ws.Cells.Item(204, 6).Value = 'Wisla Krakow II'
$ws.Cells.Item(204, 7).Value = 1
$ws.Cells.Item(204, 10).Value = 1
$ws.Cells.Item(204, 11).Value = 'A'
$ws.Cells.Item(204, 12).Value = 3.25
$ws.Cells.Item(204, 13).Value = 4
$ws.Cells.Item(204, 14).Value = 1.8
$ws.Cells.Item(204, 15).Value = 3.9
$ws.Cells.Item(204, 16).Value = 4.75
$ws.Cells.Item(204, 17).Value = 1.55
$ws.Cells.Item(204, 18).Value = 1
$ws.Cells.Item(204, 19).Value = 1.925
$ws.Cells.Item(204, 20).Value = 1.875
$ws.Cells.Item(204, 22).Value = 1.85
$ws.Cells.Item(204, 23).Value = 1.95
$ws.Cells.Item(204, 25).Value = -1
$ws.Cells.Item(204, 26).Value = 0.55
$ws.Cells.Item(204, 27).Value = 0
$ws.Cells.Item(204, 28).Value = 0
$ws.Cells.Item(204, 29).Value = -1
$ws.Cells.Item(204, 30).Value = 0.95
# Row 211
$ws.Cells.Item(211, 5).Value = 'MGKS Moto Jelcz Olawa'
